$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "26.517.08"
$ws.Cells.Item(2,5).Value = "  +6.83%  "
$ws.Cells.Item(3,4).Value = "1.718.63"
$ws.Cells.Item(3,5).Value = "  +3.36%  "
$ws.Cells.Item(4,4).Value = "'1.002"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(5,4).Value = "'333.51"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  +1.46%  "
$ws.Cells.Item(6,4).Value = "'1.001"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -0.03%  "
$ws.Cells.Item(8,4).Value = "'48.19"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "  +1.97%  "
$ws.Cells.Item(9,4).Value = "'0.3351"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "  +2.51%  "
$ws.Cells.Item(10,4).Value = "'1.182"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "  +3.81%  "
$ws.Cells.Item(11,4).Value = "'0.07369"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  +4.03%  "
$ws.Cells.Item(12,4).Value = "'1.003"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "  -0.01%  "
$ws.Cells.Item(13,4).Value = "'6.366"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "  +4.61%  "
$ws.Cells.Item(14,4).Value = "'20.03"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  +2.60%  "
$ws.Cells.Item(15,4).Value = "'7.010"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "  +5.95%  "
$ws.Cells.Item(16,4).Value = "1.717.09"
$ws.Cells.Item(16,5).Value = "  +3.36%  "
$ws.Cells.Item(17,4).Value = "'0.00001066"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "  +1.75%  "
$ws.Cells.Item(18,4).Value = "'0.06619"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -0.68%  "
$ws.Cells.Item(19,4).Value = "'81.83"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  +4.12%  "
$ws.Cells.Item(20,4).Value = "'1.001"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -0.06%  "
$ws.Cells.Item(21,4).Value = "'16.48"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +4.16%  "
$ws.Cells.Item(22,4).Value = "'6.102"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  +2.73%  "
$ws.Cells.Item(23,4).Value = "'12.72"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "  +1.14%  "
$ws.Cells.Item(24,4).Value = "26.453.99"
$ws.Cells.Item(24,5).Value = "  +6.57%  "
$ws.Cells.Item(25,4).Value = "'2.432"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  -1.56%  "
$ws.Cells.Item(26,4).Value = "'2.377"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "  -2.53%  "
$ws.Cells.Item(27,4).Value = "'1.384"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  +18.41%  "
$ws.Cells.Item(28,4).Value = "'151.53"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "  +1.14%  "
$ws.Cells.Item(30,4).Value = "1.909.82"
$ws.Cells.Item(30,5).Value = "  +3.55%  "
$ws.Cells.Item(31,4).Value = "'130.64"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = "  +3.75%  "
$ws.Cells.Item(32,4).Value = "'4.116"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = "  +1.05%  "
$ws.Cells.Item(33,4).Value = "'5.885"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  +3.28%  "
$ws.Cells.Item(34,4).Value = "'0.08605"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  +1.53%  "
$ws.Cells.Item(35,5).Value = "  +3.05%  "
$ws.Cells.Item(36,4).Value = "'12.58"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  +3.22%  "
$ws.Cells.Item(37,5).Value = "  +2.94%  "
$ws.Cells.Item(38,4).Value = "'0.02313"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  +1.73%  "
$ws.Cells.Item(39,4).Value = "'0.2150"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  +3.16%  "
$ws.Cells.Item(40,4).Value = "'0.06179"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -0.78%  "
$ws.Cells.Item(41,4).Value = "'8.395"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "  +1.77%  "
$ws.Cells.Item(42,4).Value = "'1.219"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "  -0.71%  "
$ws.Cells.Item(43,4).Value = "'0.6149"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  +3.53%  "
$ws.Cells.Item(44,4).Value = "'1.001"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "  -0.07%  "
$ws.Cells.Item(45,4).Value = "'14.08"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  +4.04%  "
$ws.Cells.Item(46,4).Value = "'3.891"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  +1.05%  "
$ws.Cells.Item(47,4).Value = "'0.5945"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "  +5.10%  "
$ws.Cells.Item(48,4).Value = "'127.97"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +1.92%  "
$ws.Cells.Item(49,4).Value = "'2.031"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "  +3.75%  "
$ws.Cells.Item(50,4).Value = "'0.07152"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +2.43%  "
$ws.Cells.Item(51,4).Value = "'76.56"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  +1.71%  "
